$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect so values/text can be updated, then
# restore protection afterwards.
$ws.Unprotect()

# Bump the "as of" date in the confidential disclaimer (row 10, column A)
# from 2021-03-22 to 2021-03-23.
[void]$ws.Cells.Replace("2021-03-22", "2021-03-23")

# Refresh the Weight (D) and Percent Change (E) figures for rows 2-7.
$ws.Cells.Item(2, 4).Value = 0.4908992694130402
$ws.Cells.Item(2, 5).Value = 0.0003954132068011429

$ws.Cells.Item(3, 4).Value = 0.3310972596284619
$ws.Cells.Item(3, 5).Value = -0.01022864019253888

$ws.Cells.Item(4, 4).Value = 0.09343237952211109
$ws.Cells.Item(4, 5).Value = -0.001339456563337227

$ws.Cells.Item(5, 4).Value = 0.05514839688456574
$ws.Cells.Item(5, 5).Value = 0.001042028482111768

$ws.Cells.Item(6, 4).Value = 0.0294226945518211
$ws.Cells.Item(6, 5).Value = -0.04427083333333326

$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = -0.004562816303544048

$ws.Protect()
